# Belgium First Division B - league bases update (02-05-2024 20:28)
#
# The source data rows got re-sorted/re-matched against a fixture list:
# for a handful of rows the "id"/result/odds data (columns B through AB)
# moved to a different row while the sequential rank in column A (and
# the shared Div/Date in C/D) stayed put. Column A is never touched here;
# only the payload columns B:AB swap/rotate between rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowPayload($rowA, $rowB) {
    $rangeA = $ws.Range("B$rowA`:AB$rowA")
    $rangeB = $ws.Range("B$rowB`:AB$rowB")
    $valA = $rangeA.Value()
    $valB = $rangeB.Value()
    $rangeA.Value = $valB
    $rangeB.Value = $valA
}

# Simple pairwise swaps
Swap-RowPayload 117 118
Swap-RowPayload 130 131
Swap-RowPayload 133 134
Swap-RowPayload 218 219
Swap-RowPayload 221 222

# 4-way rotation across rows 235-238:
#   new235 = old236, new236 = old237, new237 = old238, new238 = old235
$v235 = $ws.Range("B235:AB235").Value()
$v236 = $ws.Range("B236:AB236").Value()
$v237 = $ws.Range("B237:AB237").Value()
$v238 = $ws.Range("B238:AB238").Value()

$ws.Range("B235:AB235").Value = $v236
$ws.Range("B236:AB236").Value = $v237
$ws.Range("B237:AB237").Value = $v238
$ws.Range("B238:AB238").Value = $v235
